# "problem 40, 47, 289" — mark these LeetCode problems as solved with their
# per-category answers, on the March ("problem 46" / "problem 47" rows) and
# May ("problem 289" / "problem 39 combination-sum" / "problem 47
# permutations-ii" / "problem 40 combination-sum-ii") sheets, and move the
# active selection/tab to reflect where the author was last working (May).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# March sheet — rows 18 (problem 46, "3sum") and 20 (problem 47,
# "permutations-ii") get their Problem # and per-column answers filled in.
# ---------------------------------------------------------------------
$march = $wb.Worksheets.Item("March")

$march.Cells.Item(18, 2).Value = 46
$row18 = @("Medium", "Yes", "Yes", "Yes", "No", "No", "No", "Yes")
for ($i = 0; $i -lt $row18.Length; $i++) {
    $march.Cells.Item(18, 4 + $i).Value = $row18[$i]
}

$march.Cells.Item(20, 2).Value = 47
$row20 = @("Medium", "Yes", "No", "No", "No", "No", "No", "No")
for ($i = 0; $i -lt $row20.Length; $i++) {
    $march.Cells.Item(20, 4 + $i).Value = $row20[$i]
}

# ---------------------------------------------------------------------
# May sheet — rows 38-41 (problems 289, 39, 47, 40) get filled in; rows
# 39 and 41 (and the already-filled row 22) are marked with the "Good"
# (green) cell style to flag them as done.
# ---------------------------------------------------------------------
$may = $wb.Worksheets.Item("May")

# Row 22 was already filled in; only its highlight changes.
$may.Cells.Item(22, 3).Font.Color = 24832
$may.Cells.Item(22, 3).Interior.Color = 13561798

$may.Cells.Item(38, 2).Value = 289
$row38 = @("Medium", "Yes", "Yes", "No", "No", "No", "No", "No")
for ($i = 0; $i -lt $row38.Length; $i++) {
    $may.Cells.Item(38, 4 + $i).Value = $row38[$i]
}

$may.Cells.Item(39, 2).Value = 39
$row39 = @("Medium", "Yes", "No", "Yes", "No", "No", "No", "No")
for ($i = 0; $i -lt $row39.Length; $i++) {
    $may.Cells.Item(39, 4 + $i).Value = $row39[$i]
}
$may.Cells.Item(39, 3).Font.Color = 24832
$may.Cells.Item(39, 3).Interior.Color = 13561798

$may.Cells.Item(40, 2).Value = 47
$may.Cells.Item(40, 3).Value = "https://leetcode.com/problems/permutations-ii/"
$row40 = @("Medium", "Yes", "No", "No", "No", "No", "No", "No")
for ($i = 0; $i -lt $row40.Length; $i++) {
    $may.Cells.Item(40, 4 + $i).Value = $row40[$i]
}

$may.Cells.Item(41, 2).Value = 40
$may.Cells.Item(41, 3).Value = "https://leetcode.com/problems/combination-sum-ii/"
$row41 = @("Medium", "Yes", "No", "No", "No", "No", "No", "No")
for ($i = 0; $i -lt $row41.Length; $i++) {
    $may.Cells.Item(41, 4 + $i).Value = $row41[$i]
}
$may.Cells.Item(41, 3).Font.Color = 24832
$may.Cells.Item(41, 3).Interior.Color = 13561798

# ---------------------------------------------------------------------
# Selection / active-tab bookkeeping: work moved from March to May.
# ---------------------------------------------------------------------
$march.Range("C20").Select()

$april = $wb.Worksheets.Item("April")
$april.Range("C54").Select()

$may.Range("C39").Select()
$may.Activate()
